$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 2

# Row 7
$ws.Range("F7").Value = 2

# Row 14
$ws.Range("C14").Value = 2

# Row 17
$ws.Range("F17").Value = 2

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 2

# Row 25
$ws.Range("F25").Value = 2

# Row 26
$ws.Range("D26").Value = 2
$ws.Range("F26").Value = 2

# Update the frozen pane's top-left visible cell and the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 3
$ws.Range("C1").Select()
